$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing data row (MuSCs -> Calca/Calcr -> MuSCs) is being recomputed with
# a new TPM-based calculation AND a new row is inserted above it describing the
# MuSCs -> Calca/Calcr -> ECs interaction. Rather than using Range.Insert (which
# would clone the header's bold/bordered style onto the new row), write the
# shifted/recomputed row directly to row 3, then overwrite row 2 in place.

# Row 3: same sending/target clusters as before (MuSCs -> MuSCs), but refreshed
# expression values from the new TPM run.
$ws.Cells.Item(3, 1).Value = "MuSCs"
$ws.Cells.Item(3, 2).Value = "Calca"
$ws.Cells.Item(3, 3).Value = "Calcr"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.06056433333333333
$ws.Cells.Item(3, 8).Value = 0.181693
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 9.509424666666666
$ws.Cells.Item(3, 14).Value = 28.528274
$ws.Cells.Item(3, 15).Value = 0.9961228469411958
$ws.Cells.Item(3, 16).Value = 0.9961228469411958
$ws.Cells.Item(3, 17).Value = 0.5759319653202222
$ws.Cells.Item(3, 18).Value = 5.183387687881999
$ws.Cells.Item(3, 19).Value = 0.9961228469411958
$ws.Cells.Item(3, 20).Value = 0.9961228469411958

# Row 2: new MuSCs -> ECs interaction record.
$ws.Cells.Item(2, 1).Value = "MuSCs"
$ws.Cells.Item(2, 2).Value = "Calca"
$ws.Cells.Item(2, 3).Value = "Calcr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.06056433333333333
$ws.Cells.Item(2, 8).Value = 0.181693
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.037013
$ws.Cells.Item(2, 14).Value = 0.111039
$ws.Cells.Item(2, 15).Value = 0.003877153058804169
$ws.Cells.Item(2, 16).Value = 0.003877153058804169
$ws.Cells.Item(2, 17).Value = 0.002241667669666666
$ws.Cells.Item(2, 18).Value = 0.020175009027
$ws.Cells.Item(2, 19).Value = 0.003877153058804169
$ws.Cells.Item(2, 20).Value = 0.003877153058804169
